$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 2569.8462
$ws.Range("I11").Value = 2569.8462
$ws.Range("K11").Value = 2569.8462
$ws.Range("M11").Value = -2429.8462
$ws.Range("H41").Value = 595.86664
$ws.Range("J41").Value = 833
$ws.Range("L41").Value = 833
$ws.Range("N41").Value = -1713
$ws.Range("I62").Value = 20839932
$ws.Range("J62").Value = 8877.799999999999
$ws.Range("K62").Value = 20839932
$ws.Range("L62").Value = 8877.799999999999
$ws.Range("M62").Value = -20839308
$ws.Range("N62").Value = -10125.8
$ws.Range("H64").Value = 3427.7144
$ws.Range("I64").Value = 3509.3333
$ws.Range("K64").Value = 3509.3333
$ws.Range("M64").Value = -3261.3333
$ws.Range("I65").Value = 20839932
$ws.Range("J65").Value = 8877.799999999999
$ws.Range("K65").Value = 104199660
$ws.Range("L65").Value = 44389
$ws.Range("M65").Value = -104196540
$ws.Range("N65").Value = -50629
$ws.Range("H67").Value = 3427.7144
$ws.Range("I67").Value = 3509.3333
$ws.Range("K67").Value = 3509.3333
$ws.Range("M67").Value = -2651.3333
$ws.Range("H137").Value = 20339.033
$ws.Range("I137").Value = 27422.094
$ws.Range("K137").Value = 82266.28200000001
$ws.Range("M137").Value = -79716.28200000001
$ws.Range("H138").Value = 2500.0735
$ws.Range("J138").Value = 3087
$ws.Range("L138").Value = 9261
$ws.Range("N138").Value = -19541

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 12361
$ws.Range("I36").Value = 2500
$ws.Range("K36").Value = 2500
$ws.Range("M36").Value = -2154
$ws.Range("H61").Value = 3273.4375
$ws.Range("I61").Value = 3079.5454
$ws.Range("J61").Value = 3700
$ws.Range("K61").Value = 3079.5454
$ws.Range("L61").Value = 3700
$ws.Range("M61").Value = -2867.5454
$ws.Range("N61").Value = -4124
$ws.Range("H74").Value = 120022.766
$ws.Range("I74").Value = 135426
$ws.Range("J74").Value = 4498.5
$ws.Range("K74").Value = 135426
$ws.Range("L74").Value = 4498.5
$ws.Range("M74").Value = -134552
$ws.Range("N74").Value = -6246.5
$ws.Range("H77").Value = 120022.766
$ws.Range("I77").Value = 135426
$ws.Range("J77").Value = 4498.5
$ws.Range("K77").Value = 677130
$ws.Range("L77").Value = 22492.5
$ws.Range("M77").Value = -672762
$ws.Range("N77").Value = -31228.5
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H132").Value = 2430.641
$ws.Range("I132").Value = 2222.5881
$ws.Range("J132").Value = 2591.4092
$ws.Range("K132").Value = 6667.7643
$ws.Range("L132").Value = 7774.2276
$ws.Range("M132").Value = -4137.7643
$ws.Range("N132").Value = -12834.2276
$ws.Range("H136").Value = 3273.4375
$ws.Range("I136").Value = 3079.5454
$ws.Range("J136").Value = 3700
$ws.Range("K136").Value = 9238.636200000001
$ws.Range("L136").Value = 11100
$ws.Range("M136").Value = -6688.636200000001
$ws.Range("N136").Value = -16200

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4666.1304
$ws.Range("I99").Value = 4609.5557
$ws.Range("K99").Value = 4609.5557
$ws.Range("M99").Value = -3111.5557
$ws.Range("H132").Value = 114497.5
$ws.Range("J132").Value = 114497.5
$ws.Range("L132").Value = 114497.5
$ws.Range("N132").Value = -124617.5

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 933
$ws.Range("I7").Value = 949.5
$ws.Range("J7").Value = 900
$ws.Range("K7").Value = 949.5
$ws.Range("L7").Value = 900
$ws.Range("M7").Value = -836.5
$ws.Range("N7").Value = -1126
$ws.Range("H94").Value = 793.72
$ws.Range("I94").Value = 558.9
$ws.Range("J94").Value = 950.26666
$ws.Range("K94").Value = 558.9
$ws.Range("L94").Value = 950.26666
$ws.Range("M94").Value = -107.9
$ws.Range("N94").Value = -1852.26666
$ws.Range("H105").Value = 5203.5264
$ws.Range("I105").Value = 2138.5386
$ws.Range("J105").Value = 6797.32
$ws.Range("K105").Value = 2138.5386
$ws.Range("L105").Value = 6797.32
$ws.Range("M105").Value = -391.5385999999999
$ws.Range("N105").Value = -10291.32
$ws.Range("H122").Value = 2635.6191
$ws.Range("I122").Value = 2616.5
$ws.Range("J122").Value = 2750.3333
$ws.Range("K122").Value = 7849.5
$ws.Range("L122").Value = 8250.999899999999
$ws.Range("M122").Value = -5399.5
$ws.Range("N122").Value = -13150.9999
$ws.Range("H134").Value = 3561.8235
$ws.Range("I134").Value = 3490.0667
$ws.Range("K134").Value = 10470.2001
$ws.Range("M134").Value = -7935.2001

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 378.6
$ws.Range("J2").Value = 387.44446
$ws.Range("L2").Value = 2324.66676
$ws.Range("N2").Value = -2550.66676
$ws.Range("H68").Value = 27778662
$ws.Range("I68").Value = 83333336
$ws.Range("K68").Value = 250000008
$ws.Range("M68").Value = -249999197
$ws.Range("H71").Value = 27778662
$ws.Range("I71").Value = 83333336
$ws.Range("K71").Value = 750000024
$ws.Range("M71").Value = -749995968
$ws.Range("H76").Value = 214418190
$ws.Range("I76").Value = 300182720
$ws.Range("K76").Value = 900548160
$ws.Range("M76").Value = -900547777
$ws.Range("H79").Value = 214418190
$ws.Range("I79").Value = 300182720
$ws.Range("K79").Value = 900548160
$ws.Range("M79").Value = -900546834
$ws.Range("H92").Value = 249
$ws.Range("I92").Value = 249
$ws.Range("K92").Value = 747
$ws.Range("M92").Value = 501

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 12329.667
$ws.Range("J43").Value = 15990
$ws.Range("L43").Value = 15990
$ws.Range("N43").Value = -16292
$ws.Range("H46").Value = 15000
$ws.Range("I46").Value = 15000
$ws.Range("K46").Value = 15000
$ws.Range("M46").Value = -14844
$ws.Range("H57").Value = 19064.666
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H122").Value = 6426.579
$ws.Range("I122").Value = 5493.9165
$ws.Range("J122").Value = 8025.4287
$ws.Range("K122").Value = 16481.7495
$ws.Range("L122").Value = 24076.2861
$ws.Range("M122").Value = -14031.7495
$ws.Range("N122").Value = -28976.2861
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H140").Value = 67968.42999999999
$ws.Range("J140").Value = 67968.42999999999
$ws.Range("L140").Value = 67968.42999999999
$ws.Range("N140").Value = -78328.42999999999
$ws.Range("H141").Value = 69499.5
$ws.Range("J141").Value = 69499.5
$ws.Range("L141").Value = 69499.5
$ws.Range("N141").Value = -79859.5

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4989.125
$ws.Range("I61").Value = 4989.125
$ws.Range("K61").Value = 4989.125
$ws.Range("M61").Value = -4787.125
$ws.Range("H69").Value = 120000
$ws.Range("J69").Value = 120000
$ws.Range("L69").Value = 120000
$ws.Range("N69").Value = -121622
$ws.Range("H72").Value = 120000
$ws.Range("J72").Value = 120000
$ws.Range("L72").Value = 360000
$ws.Range("N72").Value = -368112
$ws.Range("H113").Value = 4989.125
$ws.Range("I113").Value = 4989.125
$ws.Range("K113").Value = 4989.125
$ws.Range("M113").Value = -2819.125
$ws.Range("H120").Value = 40000
$ws.Range("J120").Value = 40000
$ws.Range("L120").Value = 40000
$ws.Range("N120").Value = -49676
$ws.Range("H136").Value = 3078.8235
$ws.Range("I136").Value = 2170.3333
$ws.Range("K136").Value = 6510.999899999999
$ws.Range("M136").Value = -3960.999899999999

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 114351.89
$ws.Range("I4").Value = 145723.72
$ws.Range("J4").Value = 4550.5
$ws.Range("K4").Value = 145723.72
$ws.Range("L4").Value = 4550.5
$ws.Range("M4").Value = -145610.72
$ws.Range("N4").Value = -4776.5
$ws.Range("H82").Value = 55000
$ws.Range("J82").Value = 55000
$ws.Range("L82").Value = 55000
$ws.Range("N82").Value = -55766
$ws.Range("H85").Value = 55000
$ws.Range("J85").Value = 55000
$ws.Range("L85").Value = 55000
$ws.Range("N85").Value = -57652
$ws.Range("H136").Value = 478548.53
$ws.Range("I136").Value = 502275.94
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 1506827.82
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -1504277.82
$ws.Range("N136").Value = -17100
$ws.Range("H140").Value = 59571.5
$ws.Range("J140").Value = 59571.5
$ws.Range("L140").Value = 59571.5
$ws.Range("N140").Value = -69931.5
$ws.Range("H141").Value = 67428.625
$ws.Range("J141").Value = 67428.625
$ws.Range("L141").Value = 67428.625
$ws.Range("N141").Value = -77788.625
